# Update "Capital Eqpt Sales Tax Rate" workbook:
# 1. CESTR sheet: B1 label gets units suffix, B2 value updated to 6.8%
# 2. About sheet: remove stray note about Hong Kong workshop override
# 3. Make "About" the active/selected sheet (was "CESTR")

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsCestr = $wb.Worksheets.Item("CESTR")

# --- CESTR sheet updates ---
$wsCestr.Range("B1").Value = "Tax Rate (dimensionless)"
$wsCestr.Range("B2").Value = 0.068

# --- About sheet updates ---
# Remove the stray red note cell in B13 (and its row)
$wsAbout.Range("B13").Clear()

# --- Active sheet / selection ---
$wsAbout.Activate()
[void]$wsAbout.Range("A1").Select()
